$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "account statement" detail table (rows 16-21) is re-sorted from
# grouped-by-worker / descending period (1903,1902,1901) into
# grouped-by-period / (1901,1901,1902,1902,1903,1903) with both workers'
# records interleaved. At the same time the "Salario Basico" (column G)
# for JUAN DARIO LOMBANA HERRERA is corrected from 781242 to 1700000
# (now matching DAVID GUILLERMO RONDON VISBAL), and his "Valor Mora"
# (column F) for period 1901 is corrected from 56667 to 68000.

# Row 16: JUAN DARIO LOMBANA HERRERA - periodo 1901
$ws.Range("C16").Value = "1051818784"
$ws.Range("D16").Value = "JUAN DARIO LOMBANA HERRERA"
$ws.Range("E16").Value = "1901"
$ws.Range("F16").Value = 68000
$ws.Range("G16").Value = 1700000

# Row 17: DAVID GUILLERMO RONDON VISBAL - periodo 1901
$ws.Range("C17").Value = "84458883"
$ws.Range("D17").Value = "DAVID GUILLERMO RONDON VISBAL"
$ws.Range("E17").Value = "1901"
$ws.Range("F17").Value = 68000
$ws.Range("G17").Value = 1700000

# Row 18: JUAN DARIO LOMBANA HERRERA - periodo 1902
$ws.Range("C18").Value = "1051818784"
$ws.Range("D18").Value = "JUAN DARIO LOMBANA HERRERA"
$ws.Range("E18").Value = "1902"
$ws.Range("F18").Value = 68000
$ws.Range("G18").Value = 1700000

# Row 19: DAVID GUILLERMO RONDON VISBAL - periodo 1902
$ws.Range("C19").Value = "84458883"
$ws.Range("D19").Value = "DAVID GUILLERMO RONDON VISBAL"
$ws.Range("E19").Value = "1902"
$ws.Range("F19").Value = 68000
$ws.Range("G19").Value = 1700000

# Row 20: JUAN DARIO LOMBANA HERRERA - periodo 1903
$ws.Range("C20").Value = "1051818784"
$ws.Range("D20").Value = "JUAN DARIO LOMBANA HERRERA"
$ws.Range("E20").Value = "1903"
$ws.Range("F20").Value = 56667
$ws.Range("G20").Value = 1700000

# Row 21: DAVID GUILLERMO RONDON VISBAL - periodo 1903
$ws.Range("C21").Value = "84458883"
$ws.Range("D21").Value = "DAVID GUILLERMO RONDON VISBAL"
$ws.Range("E21").Value = "1903"
$ws.Range("F21").Value = 56667
$ws.Range("G21").Value = 1700000

$wb.Save()
